$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.310.80"
$ws.Range("E2").Value = "  +0.58%  "

$ws.Range("D3").Value = "1.865.96"
$ws.Range("E3").Value = "  +0.43%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.51%  "

$ws.Range("E6").Value = "  +0.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4675"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.41%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2837"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.54%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06519"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.53%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.81"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.00%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07934"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.00%  "

$ws.Range("D13").Value = "1.871.24"
$ws.Range("E13").Value = "  +0.70%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.150"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.14%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6778"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.93%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "279.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.09%  "

$ws.Range("D17").Value = "30.303.29"
$ws.Range("E17").Value = "  +0.44%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.54%  "

$ws.Range("E19").Value = "  +0.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.410"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.93%  "

$ws.Range("D21").Value = "2.116.47"
$ws.Range("E21").Value = "  +0.93%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.000007303"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.74%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.151"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "166.43"
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.157"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.68%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.09"
$ws.Range("D27").Style = "Normal"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.931"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.03%  "

$ws.Range("E29").Value = "  +3.40%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09724"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.87%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.398"
$ws.Range("D31").Style = "Normal"

$ws.Range("E32").Value = "  +0.35%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.083"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.66%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04747"
$ws.Range("D34").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.131"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.95%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7075"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.39%  "

$ws.Range("E37").Value = "  +0.29%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01864"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.02%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.574"
$ws.Range("D39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.314"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.59"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.46%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.968"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.26%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8500"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.40%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4190"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.61%  "

$ws.Range("E45").Value = "  +0.09%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "967.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.78%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.187"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.73%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.328"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.97%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.78%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1131"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.18%  "
